$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"

# --- Narrow the status columns ---
$overview.Columns.Item(5).ColumnWidth = 13.4101845877511
$overview.Columns.Item(6).ColumnWidth = 13.4101845877511

$zhcn.Columns.Item(3).ColumnWidth = 13.4101845877511

$dede.Columns.Item(3).ColumnWidth = 13.4101845877511
